$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = "Автонеева";  C = "Катерина"; D = "Николаевна" },
    @{ Row = 3;  B = "Щиборин";    C = "Владимир"; D = "Георгиевич" },
    @{ Row = 4;  B = "Сивриков";   C = "Дмитрий";  D = "Петрович" },
    @{ Row = 5;  B = "Позов";      C = "Леонид";   D = "Николаевич" },
    @{ Row = 6;  B = "Сымшикова";  C = "Алла";     D = "Василевна" },
    @{ Row = 7;  B = "Шогина";     C = "Полина";   D = "Евгеньевна" },
    @{ Row = 8;  B = "Налютин";    C = "Олег";     D = "Евгеньевич" },
    @{ Row = 9;  B = "Буртынкин";  C = "Вячеслав"; D = "Петрович" },
    @{ Row = 10; B = "Пикарова";   C = "Ольга";    D = "Олеговна" },
    @{ Row = 11; B = "Годовкин";   C = "Илья";     D = "Юриевич" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
